$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(68, 1).Value = "2025-04-29 09:50:33"
$ws.Cells.Item(68, 2).Value = 200
